$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-16 04:48:31'
$ws.Range('E3').Value = '2026-02-16 04:48:33'
$ws.Range('E4').Value = '2026-02-16 04:48:35'
$ws.Range('J4').Value = '1014.5 hPa'
$ws.Range('K4').Value = '-0.1 MJ/m2'
$ws.Range('E5').Value = '2026-02-16 04:48:38'
$ws.Range('M5').Value = '-1.0 °C 4:29 TU'
$ws.Range('N5').Value = '-1.2 °C 4:07 TU'
$ws.Range('E6').Value = '2026-02-16 04:48:40'
$ws.Range('J6').Value = '1014.6 hPa'
$ws.Range('N6').Value = '5.9 °C 4:11 TU'
$ws.Range('E7').Value = '2026-02-16 04:48:43'
$ws.Range('H7').NumberFormat = "@"
$ws.Range('H7').Value = '55%'
$ws.Range('J7').Value = '1014.9 hPa'
$ws.Range('M7').Value = '13.5 °C 4:20 TU'
$ws.Range('E8').Value = '2026-02-16 04:48:45'
$ws.Range('J8').Value = '1014.7 hPa'
$ws.Range('M8').Value = '9.7 °C 4:21 TU'
$ws.Range('E9').Value = '2026-02-16 04:48:47'
$ws.Range('O9').Value = '5.2 °C'
$ws.Range('E10').Value = '2026-02-16 04:48:50'
$ws.Range('M10').Value = '4.6 °C 4:08 TU'
$ws.Range('O10').Value = '4.0 °C'
$ws.Range('E11').Value = '2026-02-16 04:48:52'
$ws.Range('N11').Value = '0.3 °C 4:29 TU'
$ws.Range('E12').Value = '2026-02-16 04:48:54'
$ws.Range('N12').Value = '4.8 °C 4:12 TU'
$ws.Range('E13').Value = '2026-02-16 04:48:57'
$ws.Range('N13').Value = '0.1 °C 4:19 TU'
$ws.Range('O13').Value = '1.2 °C'
$ws.Range('E14').Value = '2026-02-16 04:48:59'
$ws.Range('E15').Value = '2026-02-16 04:49:01'
$ws.Range('H15').NumberFormat = "@"
$ws.Range('H15').Value = '90%'
$ws.Range('N15').Value = '3.7 °C 4:29 TU'
$ws.Range('O15').Value = '5.4 °C'
$ws.Range('E16').Value = '2026-02-16 04:49:04'
$ws.Range('H16').NumberFormat = "@"
$ws.Range('H16').Value = '78%'
$ws.Range('M16').Value = '0.5 °C 4:20 TU'
$ws.Range('O16').Value = '-0.7 °C'
$ws.Range('E17').Value = '2026-02-16 04:49:06'
$ws.Range('N17').Value = '4.9 °C 4:19 TU'
$ws.Range('E18').Value = '2026-02-16 04:49:08'
$ws.Range('J18').Value = '1015.0 hPa'
$ws.Range('O18').Value = '4.3 °C'
$ws.Range('E19').Value = '2026-02-16 04:49:11'
$ws.Range('N19').Value = '2.6 °C 4:06 TU'
$ws.Range('O19').Value = '3.3 °C'
$ws.Range('E20').Value = '2026-02-16 04:49:13'
$ws.Range('H20').NumberFormat = "@"
$ws.Range('H20').Value = '91%'
$ws.Range('N20').Value = '-2.0 °C 4:17 TU'
$ws.Range('O20').Value = '-1.2 °C'
$ws.Range('E21').Value = '2026-02-16 04:49:15'
$ws.Range('H21').NumberFormat = "@"
$ws.Range('H21').Value = '81%'
$ws.Range('N21').Value = '3.3 °C 4:15 TU'
$ws.Range('O21').Value = '4.8 °C'
$ws.Range('E22').Value = '2026-02-16 04:49:18'
$ws.Range('I22').Value = '0.6 mm'
$ws.Range('N22').Value = '-6.5 °C 4:08 TU'
$ws.Range('E23').Value = '2026-02-16 04:49:20'
$ws.Range('L23').Value = '49.7 km/h - 325º 4:27 TU'
$ws.Range('E24').Value = '2026-02-16 04:49:22'
$ws.Range('H24').NumberFormat = "@"
$ws.Range('H24').Value = '71%'
$ws.Range('J24').Value = '1018.1 hPa'
$ws.Range('E25').Value = '2026-02-16 04:49:25'
$ws.Range('O25').Value = '0.6 °C'
$ws.Range('E26').Value = '2026-02-16 04:49:27'
$ws.Range('E27').Value = '2026-02-16 04:49:29'
$ws.Range('H27').NumberFormat = "@"
$ws.Range('H27').Value = '75%'
$ws.Range('N27').Value = '0.1 °C 4:29 TU'
$ws.Range('O27').Value = '0.9 °C'
$ws.Range('E28').Value = '2026-02-16 04:49:32'
$ws.Range('J28').Value = '1015.9 hPa'
$ws.Range('N28').Value = '2.5 °C 4:29 TU'
$ws.Range('E29').Value = '2026-02-16 04:49:34'
$ws.Range('E30').Value = '2026-02-16 04:49:36'
$ws.Range('J30').Value = '1014.6 hPa'
$ws.Range('E31').Value = '2026-02-16 04:49:39'
$ws.Range('J31').Value = '1013.1 hPa'
$ws.Range('E32').Value = '2026-02-16 04:49:41'
$ws.Range('E33').Value = '2026-02-16 04:49:44'
$ws.Range('H33').NumberFormat = "@"
$ws.Range('H33').Value = '72%'
$ws.Range('J33').Value = '1016.0 hPa'
$ws.Range('N33').Value = '2.4 °C 4:29 TU'
$ws.Range('O33').Value = '4.6 °C'
$ws.Range('E34').Value = '2026-02-16 04:49:46'
$ws.Range('E35').Value = '2026-02-16 04:49:48'
$ws.Range('N35').Value = '6.6 °C 4:16 TU'
$ws.Range('O35').Value = '6.7 °C'
$ws.Range('E36').Value = '2026-02-16 04:49:51'
$ws.Range('H36').NumberFormat = "@"
$ws.Range('H36').Value = '91%'
$ws.Range('J36').Value = '1014.5 hPa'
$ws.Range('N36').Value = '5.2 °C 4:05 TU'
$ws.Range('O36').Value = '6.8 °C'
$ws.Range('E37').Value = '2026-02-16 04:49:53'
$ws.Range('H37').NumberFormat = "@"
$ws.Range('H37').Value = '94%'
$ws.Range('J37').Value = '1018.2 hPa'
$ws.Range('N37').Value = '0.7 °C 4:29 TU'
$ws.Range('O37').Value = '1.8 °C'
$ws.Range('E38').Value = '2026-02-16 04:49:56'
$ws.Range('M38').Value = '6.8 °C 4:18 TU'
$ws.Range('E39').Value = '2026-02-16 04:49:58'
$ws.Range('M39').Value = '0.7 °C 4:16 TU'
$ws.Range('E40').Value = '2026-02-16 04:50:01'
$ws.Range('J40').Value = '1019.3 hPa'
$ws.Range('N40').Value = '1.5 °C 4:17 TU'
$ws.Range('O40').Value = '3.0 °C'
$ws.Range('E41').Value = '2026-02-16 04:50:03'
$ws.Range('J41').Value = '1016.0 hPa'
$ws.Range('O41').Value = '14.9 °C'
$ws.Range('E42').Value = '2026-02-16 04:50:05'
$ws.Range('H42').NumberFormat = "@"
$ws.Range('H42').Value = '96%'
$ws.Range('N42').Value = '5.7 °C 4:27 TU'
$ws.Range('E43').Value = '2026-02-16 04:50:08'
$ws.Range('E44').Value = '2026-02-16 04:50:10'
$ws.Range('L44').Value = '44.3 km/h - 51º 4:24 TU'
$ws.Range('E45').Value = '2026-02-16 04:50:13'
$ws.Range('J45').Value = '1019.7 hPa'
$ws.Range('E46').Value = '2026-02-16 04:50:15'
$ws.Range('J46').Value = '1018.6 hPa'
$ws.Range('L46').Value = '38.5 km/h - 328º 4:13 TU'
$ws.Range('O46').Value = '12.6 °C'
